## "entry overrides for * and -"
## Adds a new slide 2, duplicated from slide 1 ("Tlön: the Beginning" /
## class Foo.Bar(name:str){ add := ... }) and retargets the copy's text to
## the new "Tuples and Class Functions" / class My.Point(xy:(i32,i32))
## content, matching the author's edit.

$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

# Duplicate slide 1 -> lands immediately after it as the new slide 2.
# This preserves shape ids/names/xfrm and the p14:creationId ext exactly
# as a real PowerPoint "duplicate slide" would, which is what the target
# slide2.xml looks like (same shape tree shape as slide1, new text).
$range = $s1.Duplicate()
$s2 = $p.Slides.Item(2)

# --- Body placeholder: "class Foo.Bar(name:str){ ... }" -----------------
$body = $s2.Shapes.Item(1)
$tr = $body.TextFrame.TextRange

# Work back-to-front so earlier character offsets stay valid.
# Layout of the original text (1-based TextRange offsets):
#   [1-6]   "class "
#   [7-13]  "Foo.Bar"
#   [14]    "("
#   [15-22] "name:str"
#   [23]    ")"
#   [24]    <br>
#   [25]    "{"
#   [26]    <br>
#   [27-52] "  add := (x,y:i32) -> i32;"
#   [53]    <br>
#   [54]    "}"

# "  add := (x,y:i32) -> i32;" -> "  " + "distance := () -> i32 { … }"
# (leave the leading two spaces alone so the run splits in two, as in
# the target XML)
$tr.Characters(29, 24).Text = "distance := () -> i32 { … }"

# ")" -> ":(i32,i32))"
$tr.Characters(23, 1).Text = ":(i32,i32))"

# "name:str" -> "xy"
$tr.Characters(15, 8).Text = "xy"

# "Foo.Bar" -> "My.Point"
$tr.Characters(7, 7).Text = "My.Point"

# --- Title placeholder: "Tlön: the Beginning" ---------------------------
$title = $s2.Shapes.Item(2)
$ttr = $title.TextFrame.TextRange

# Drop the first run ("Tlön") entirely, then retext the remaining run
# (": the Beginning") to the new title, which keeps that run's rPr
# (no spelling-error "err" flag) exactly as the target does.
$ttr.Characters(1, 4).Text = ""
$ttr2 = $title.TextFrame.TextRange
$ttr2.Characters(1, $ttr2.Length).Text = "Tuples and Class Functions"
